# Update cell values per the target diff (row-by-row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 26.3
$ws.Range("D4").Value = 543
$ws.Range("F4").Value = 0.2
$ws.Range("N4").Value = 21.3
# Row 7
$ws.Range("D7").Value = 7090
$ws.Range("F7").Value = 30.2
$ws.Range("N7").Value = 4.8
# Row 8
$ws.Range("D8").Value = 12.4
# Row 11
$ws.Range("C11").Value = 9.199999999999999
$ws.Range("D11").Value = 145.4
$ws.Range("F11").Value = 0.2
# Row 12
$ws.Range("C12").Value = 1.6
$ws.Range("D12").Value = 25
$ws.Range("N12").Value = 0.9
# Row 15
$ws.Range("C15").Value = 31
$ws.Range("D15").Value = 376.5
# Row 16
$ws.Range("C16").Value = 0.7
$ws.Range("D16").Value = 16.8
# Row 19
$ws.Range("C19").Value = 10.8
$ws.Range("D19").Value = 233.4
# Row 20
$ws.Range("C20").Value = 3.4
$ws.Range("D20").Value = 89.59999999999999
# Row 23
$ws.Range("C23").Value = 73.5
$ws.Range("D23").Value = 1555.5
$ws.Range("F23").Value = 8.4
# Row 24
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 24.8
$ws.Range("N24").Value = 1.1
# Row 27
$ws.Range("C27").Value = 24.4
$ws.Range("D27").Value = 356.4
# Row 28
$ws.Range("C28").Value = 5.8
$ws.Range("D28").Value = 102.5
$ws.Range("F28").Value = 0.1
$ws.Range("N28").Value = 3.3
# Row 31
$ws.Range("D31").Value = 1445.3
# Row 32
$ws.Range("C32").Value = 0.9
$ws.Range("D32").Value = 30.1
$ws.Range("N32").Value = 0.9
# Row 35
$ws.Range("C35").Value = 33.8
$ws.Range("D35").Value = 506.4
# Row 36
$ws.Range("C36").Value = 0.4
$ws.Range("D36").Value = 12.1
$ws.Range("F36").Value = 0.1
# Row 39
$ws.Range("D39").Value = 182.2
# Row 40
$ws.Range("C40").Value = 0.6
$ws.Range("D40").Value = 22
$ws.Range("N40").Value = 0.3
# Row 43
$ws.Range("C43").Value = 14.6
$ws.Range("D43").Value = 263.2
$ws.Range("F43").Value = 0.9
# Row 44
$ws.Range("C44").Value = 0.4
$ws.Range("D44").Value = 9.4
$ws.Range("N44").Value = 0.8
# Row 47
$ws.Range("C47").Value = 10.5
$ws.Range("D47").Value = 138
$ws.Range("F47").Value = 0.5
# Row 48
$ws.Range("C48").Value = 0.4
$ws.Range("D48").Value = 12.7
# Row 51
$ws.Range("C51").Value = 5
$ws.Range("D51").Value = 143.3
$ws.Range("E51").ClearContents()
# Row 52
$ws.Range("C52").Value = 2.8
$ws.Range("D52").Value = 68.5
$ws.Range("N52").Value = 2.2
# Row 55
$ws.Range("C55").Value = 40.2
$ws.Range("D55").Value = 714
$ws.Range("E55").Value = 0.9
$ws.Range("F55").Value = 3.2
# Row 56
$ws.Range("C56").Value = 0.7
$ws.Range("D56").Value = 24.6
$ws.Range("N56").Value = 0.9
# Row 59
$ws.Range("C59").Value = 28.2
$ws.Range("D59").Value = 386.1
$ws.Range("F59").Value = 1.6
# Row 60
$ws.Range("C60").Value = 0.4
$ws.Range("D60").Value = 19.8
$ws.Range("N60").Value = 0.9
# Row 63
$ws.Range("C63").Value = 10.2
$ws.Range("D63").Value = 236
# Row 64
$ws.Range("C64").Value = 1.2
$ws.Range("D64").Value = 24.9
$ws.Range("N64").Value = 1.4
# Row 67
$ws.Range("C67").Value = 13.3
$ws.Range("D67").Value = 282.6
$ws.Range("F67").Value = 1
# Row 68
$ws.Range("C68").Value = 0.9
$ws.Range("D68").Value = 36.1
$ws.Range("N68").Value = 0.6
# Row 71
$ws.Range("C71").Value = 25.8
$ws.Range("D71").Value = 449.1
# Row 72
$ws.Range("C72").Value = 0.8
$ws.Range("D72").Value = 22.7
$ws.Range("N72").Value = 1.3
# Row 75
$ws.Range("C75").Value = 15.9
$ws.Range("D75").Value = 258.3
$ws.Range("F75").Value = 0.8
# Row 76
$ws.Range("C76").Value = 1
$ws.Range("D76").Value = 16.4
# Row 79
$ws.Range("C79").Value = 10.4
$ws.Range("D79").Value = 179.4
$ws.Range("F79").Value = 0.7
# Row 80
$ws.Range("D80").Value = 10.2
# Row 83
$ws.Range("D83").Value = 126.2
# Row 84
$ws.Range("C84").Value = 3.9
$ws.Range("D84").Value = 107.7
$ws.Range("F84").ClearContents()
$ws.Range("N84").Value = 6.5
# Row 87
$ws.Range("D87").Value = 1400.1
$ws.Range("F87").Value = 4.7
# Row 88
$ws.Range("C88").Value = 1.1
$ws.Range("D88").Value = 23.2
# Row 91
$ws.Range("C91").Value = 21.1
$ws.Range("D91").Value = 296.9
$ws.Range("F91").Value = 0.7
# Row 96
$ws.Range("C96").Value = 1.9
$ws.Range("D96").Value = 55.3
$ws.Range("N96").Value = 1
# Row 99
$ws.Range("D99").Value = 540.4
$ws.Range("F99").Value = 1
